$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column for all data rows (2-432)
# from serial date 45204 (2023-10-05) to 45205 (2023-10-06).
$ws.Range("C2:C432").Value = 45205
